$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Done"
$ws.Range("D3").Value = "In Progress"

$ws.Range("D10").Select()
